# Petty cash book - 11-Feb-2021 midday update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 11: Wages Expense on 44236 - add an extra 300000 ---
$ws.Range("D11").Formula = "=60000+300000"

# --- Row 12: TRANSFER BCA - correct 842000 to 842500, add two more payments ---
$ws.Range("D12").Formula = "=8839000+6709000+3720000+842500+25000000+900000+1488000"

# --- Row 13: A/R - add two more receipts ---
$ws.Range("C13").Formula = "=6000000+4332500+56547500"

# --- Row 14: new entry - SALES - cash/retail ---
$ws.Range("B14").Value = "SALES - cash/retail"
$ws.Range("C14").Formula = "=36539975+27017525-56547500"

# --- Row 15: new entry - SELISIH - kurang ---
$ws.Range("B15").Value = "SELISIH - kurang"
$ws.Range("D15").Value = 70500

# --- Row 16: new entry - SETOR KE BANK ---
$ws.Range("B16").Value = "SETOR KE BANK"
$ws.Range("D16").Value = 26000000

# --- Row 17: new day 44237 - Wages Expense ---
$ws.Range("A17").Value = 44237
$ws.Range("B17").Value = "Wages Expense"
$ws.Range("D17").Formula = "=60000+260000"

# --- Row 18: A/R ---
$ws.Range("B18").Value = "A/R"
$ws.Range("C18").Formula = "=6550000+10079000"

# --- Row 19: FREIGHT OUT ---
$ws.Range("B19").Value = "FREIGHT OUT"
$ws.Range("D19").Value = 135000

# --- Row 20: TRANSFER BCA ---
$ws.Range("B20").Value = "TRANSFER BCA"
$ws.Range("D20").Formula = "=5175000+3900000+1004000+70000"

# --- Row 21: SALES - cash/retail ---
$ws.Range("B21").Value = "SALES - cash/retail"
$ws.Range("C21").Formula = "=3146975+11687025-10079000"

# --- Row 22: SELISIH - lebih ---
$ws.Range("B22").Value = "SELISIH - lebih"
$ws.Range("C22").Value = 114500

# --- Row 23: SETOR KE BANK ---
$ws.Range("B23").Value = "SETOR KE BANK"
$ws.Range("D23").Value = 11000000

# --- Row 24: new day 44238 - Wages Expense ---
$ws.Range("A24").Value = 44238
$ws.Range("B24").Value = "Wages Expense"
$ws.Range("D24").Formula = "=60000"

# --- Row 25: A/R ---
$ws.Range("B25").Value = "A/R"
$ws.Range("C25").Formula = "=2000000+762500+5840000"

# --- Row 26: TRANSFER BCA ---
$ws.Range("B26").Value = "TRANSFER BCA"
$ws.Range("D26").Formula = "=762500"

# --- Row 27: CHEQUE RECEIVED ---
$ws.Range("B27").Value = "CHEQUE RECEIVED"
$ws.Range("D27").Formula = "=3155000"

# --- Selection state: cursor left on E44 after the midday update ---
$ws.Activate()
$ws.Range("E44").Select()
